# Insert two new data rows at the top of the "Albahaca" price-record block
# (old row 277 onward shift down by 2, becoming rows 279..319).
# New rows 277 and 278 receive their own fresh record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 277 (pushes old rows 277-317
# down to become rows 279-319), copying formatting from the row above as
# Excel normally does for EntireRow.Insert().
$ws.Rows.Item(277).EntireRow.Insert()
$ws.Rows.Item(278).EntireRow.Insert()

# --- New row 277 ---
$ws.Cells.Item(277, 1).Value = 10
$ws.Cells.Item(277, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(277, 3).Value = "La Araucanía"
$ws.Cells.Item(277, 4).Value = 44946
$ws.Cells.Item(277, 5).Value = 9
$ws.Cells.Item(277, 6).Value = 100112052
$ws.Cells.Item(277, 7).Value = "Albahaca"
$ws.Cells.Item(277, 8).Value = "Sin especificar"
$ws.Cells.Item(277, 9).Value = "Primera"
$ws.Cells.Item(277, 10).Value = 65
$ws.Cells.Item(277, 11).Value = 4000
$ws.Cells.Item(277, 12).Value = 4000
$ws.Cells.Item(277, 13).Value = 4000
$ws.Cells.Item(277, 14).Value = "$/paquete"
$ws.Cells.Item(277, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(277, 16).Value = 4000
$ws.Cells.Item(277, 17).Value = 1
$ws.Cells.Item(277, 18).Value = "Hortaliza"

# --- New row 278 ---
$ws.Cells.Item(278, 1).Value = 10
$ws.Cells.Item(278, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(278, 3).Value = "La Araucanía"
$ws.Cells.Item(278, 4).Value = 44946
$ws.Cells.Item(278, 5).Value = 9
$ws.Cells.Item(278, 6).Value = 100112052
$ws.Cells.Item(278, 7).Value = "Albahaca"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 70
$ws.Cells.Item(278, 11).Value = 4000
$ws.Cells.Item(278, 12).Value = 5000
$ws.Cells.Item(278, 13).Value = 4643
$ws.Cells.Item(278, 14).Value = "$/paquete"
$ws.Cells.Item(278, 15).Value = "Región del Maule"
$ws.Cells.Item(278, 16).Value = 4643
$ws.Cells.Item(278, 17).Value = 1
$ws.Cells.Item(278, 18).Value = "Hortaliza"
